$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-03-12 Tuesday"; new = "2024-03-13 Wednesday"},
    @{old = "108÷7=15, 3"; new = "322÷3=107, 1"},
    @{old = "337÷8=42, 1"; new = "215÷7=30, 5"},
    @{old = "774÷2=387, 0"; new = "503÷3=167, 2"},
    @{old = "978÷8=122, 2"; new = "826÷2=413, 0"},
    @{old = "541÷4=135, 1"; new = "815÷3=271, 2"},
    @{old = "453÷2=226, 1"; new = "608÷6=101, 2"},
    @{old = "860÷5=172, 0"; new = "138÷3=46, 0"},
    @{old = "898÷2=449, 0"; new = "796÷8=99, 4"},
    @{old = "750÷3=250, 0"; new = "779÷3=259, 2"},
    @{old = "588÷7=84, 0"; new = "995÷8=124, 3"},
    @{old = "765÷4=191, 1"; new = "430÷4=107, 2"},
    @{old = "388÷9=43, 1"; new = "508÷3=169, 1"},
    @{old = "112÷2=56, 0"; new = "653÷6=108, 5"},
    @{old = "584÷8=73, 0"; new = "952÷3=317, 1"},
    @{old = "957÷2=478, 1"; new = "684÷9=76, 0"},
    @{old = "534÷9=59, 3"; new = "976÷5=195, 1"},
    @{old = "621÷8=77, 5"; new = "393÷7=56, 1"},
    @{old = "671÷8=83, 7"; new = "760÷5=152, 0"},
    @{old = "892÷4=223, 0"; new = "521÷8=65, 1"},
    @{old = "158÷6=26, 2"; new = "450÷8=56, 2"},
    @{old = "325÷7=46, 3"; new = "823÷8=102, 7"},
    @{old = "373÷6=62, 1"; new = "395÷3=131, 2"},
    @{old = "307÷7=43, 6"; new = "561÷5=112, 1"},
    @{old = "401÷4=100, 1"; new = "240÷9=26, 6"},
    @{old = "731÷7=104, 3"; new = "429÷3=143, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
